# Generate Report for Handoff
# Updates the localization-status workbook to reflect that the
# a94a3e9b-c574-4087-9961-b9d0b4140a4e handback file is ready for a new
# handoff (its handback version is stale vs. the latest source).

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dea401f634fd5e86146ef630d98791c7b5014b0b/e2e/a94a3e9b-c574-4087-9961-b9d0b4140a4e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/701f5bff7b2d54db1ba02b2284f989cc9f16be0f/e2e/a94a3e9b-c574-4087-9961-b9d0b4140a4e.md."

# ---- Overview sheet: row 3 is the a94a3e9b-c574-4087-9961-b9d0b4140a4e.md file ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-21 18:52:17"

# ---- zh-cn sheet: row 3 ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-21 18:52:12"
$wsZhCn.Range("P3").Value = $errorDetail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---- de-de sheet: row 3 ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-21 18:52:17"
$wsDeDe.Range("P3").Value = $errorDetail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
